# Apply the "Handle Parallel in C#" change to the API_Controller sheet:
# add a new UserClickEvent block (IdUser/Guid, IdArticle/Guid) plus a
# merged notes cell describing the dedupe/ranking algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API_Controller")

# Header: I1:J1 merged "UserClickEvent" (same style as F1:G1 "UserInteraction")
$ws.Range("I1").Value = "UserClickEvent"
$ws.Range("I1:J1").Merge()
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4108

# Field rows, mirroring F2:G3
$ws.Range("I2").Value = "IdUser"
$ws.Range("J2").Value = "Guid"
$ws.Range("I3").Value = "IdArticle"
$ws.Range("J3").Value = "Guid"

# Notes block I5:J9 merged, centered + top aligned + wrap text
$notes = "1. Gộp nhóm lại có ArticleId cùng channel, cùng `ncategory,`n2. Đếm và sắp xếp giảm dần theo channel, cate`n3. chọn ra các bản tin tương tự khoảng 50% trong top 10`ntheo phần trăm tổng 60%`n40% còn lại  xử lý trùng các kênh đã đăng ký"
$ws.Range("I5").Value = $notes
$ws.Range("I5:J9").Merge()
$ws.Range("I5:J9").HorizontalAlignment = -4108
$ws.Range("I5:J9").VerticalAlignment = -4160
$ws.Range("I5:J9").WrapText = $true

# Column J width to fit the new notes text
$ws.Columns.Item(10).ColumnWidth = 33.6328125

# Update selection to match the author's final cursor position
$ws.Range("I5:J9").Select()
